$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Ready for handoff" report generation for file b.md
#
# The b.md row is being moved from "Handed back: in sync with en-US" to
# "Ready for handoff" on every sheet, and a brand-new handoff file/timestamp
# is recorded for it on the per-locale sheets (zh-cn, de-de). The other rows
# (a.md, .localization-config) are untouched.
# ---------------------------------------------------------------------------

$newStatus       = "Ready for handoff"
$zhHandoffFile   = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhHandoffDate   = "2016-02-24 08:53:22"
$deHandoffFile   = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$deHandoffDate   = "2016-02-24 08:53:35"

# Hyperlink colour/underline used throughout this workbook for link cells
# (matches the custom "HyperLink" cell style already present in the file).
# Excel's Font.Color takes an OLE_COLOR (0x00BBGGRR), so RGB(0x64,0x95,0xED)
# -> FF6495ED encodes as 0xED9564 == 15570276.
$linkColor = 15570276

function Set-LinkCellFormat($range) {
    $range.Font.Underline = 2
    $range.Font.Color = $linkColor
}

# ---------------------------------------------------------------------------
# Overview sheet: B3/C3 plain-text status for b.md
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B3").Value = $newStatus
$wsZh.Range("C3").Value = $zhHandoffFile
$wsZh.Range("D3").Value = $zhHandoffDate

# The hyperlink collection in this runtime can only be cleared in one shot
# (deleting a single item is a no-op), so capture every existing hyperlink's
# target/display text first, wipe them all, then recreate them - updating
# only the C3 entry to point at its new display text.
$zhLinks = @(
    @{ Cell = "A2"; Url = "https://github.com/OpenLocalizationTest/oltest/blob/5c07f89c345e59416395c218e681e08e80c6f9f4/e2e/a.md"; Text = "a.md" },
    @{ Cell = "C2"; Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b2387e355bf7069f3f3c5fb358c4909a3d9f60b7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; Text = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" },
    @{ Cell = "E2"; Url = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/3f9493c5e42ad61da462c5f4ce034efe34f90c8d/e2e/a.md"; Text = "a.md" },
    @{ Cell = "F2"; Url = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/de68fdf1899316cc2ccfa26458236f099ea823c5/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; Text = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" },
    @{ Cell = "A3"; Url = "https://github.com/OpenLocalizationTest/oltest/blob/5c07f89c345e59416395c218e681e08e80c6f9f4/e2e/b.md"; Text = "b.md" },
    @{ Cell = "C3"; Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b2387e355bf7069f3f3c5fb358c4909a3d9f60b7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; Text = $zhHandoffFile },
    @{ Cell = "E3"; Url = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/3f9493c5e42ad61da462c5f4ce034efe34f90c8d/e2e/a.md"; Text = "a.md" },
    @{ Cell = "F3"; Url = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/de68fdf1899316cc2ccfa26458236f099ea823c5/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; Text = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" },
    @{ Cell = "A4"; Url = "https://github.com/OpenLocalizationTest/oltest/blob/5c07f89c345e59416395c218e681e08e80c6f9f4/.localization-config"; Text = ".localization-config" }
)

$wsZh.Range("A1").Hyperlinks.Delete()
foreach ($link in $zhLinks) {
    $rng = $wsZh.Range($link.Cell)
    $wsZh.Hyperlinks.Add($rng, $link.Url, "", "", $link.Text)
    Set-LinkCellFormat $rng
}

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B3").Value = $newStatus
$wsDe.Range("C3").Value = $deHandoffFile
$wsDe.Range("D3").Value = $deHandoffDate

$deLinks = @(
    @{ Cell = "A2"; Url = "https://github.com/OpenLocalizationTest/oltest/blob/5c07f89c345e59416395c218e681e08e80c6f9f4/e2e/a.md"; Text = "a.md" },
    @{ Cell = "C2"; Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/245ab08e2c37ad855a4c99678cb3b7dadee35318/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; Text = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" },
    @{ Cell = "E2"; Url = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/f5b6695674ffd81eec4e46848fc0da3edbe7e98b/e2e/a.md"; Text = "a.md" },
    @{ Cell = "F2"; Url = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9bfb384904e32f7e60b4302b42b3744ef2a4b4b6/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; Text = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" },
    @{ Cell = "A3"; Url = "https://github.com/OpenLocalizationTest/oltest/blob/5c07f89c345e59416395c218e681e08e80c6f9f4/e2e/b.md"; Text = "b.md" },
    @{ Cell = "C3"; Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/245ab08e2c37ad855a4c99678cb3b7dadee35318/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; Text = $deHandoffFile },
    @{ Cell = "E3"; Url = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/f5b6695674ffd81eec4e46848fc0da3edbe7e98b/e2e/a.md"; Text = "a.md" },
    @{ Cell = "F3"; Url = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9bfb384904e32f7e60b4302b42b3744ef2a4b4b6/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; Text = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" },
    @{ Cell = "A4"; Url = "https://github.com/OpenLocalizationTest/oltest/blob/5c07f89c345e59416395c218e681e08e80c6f9f4/.localization-config"; Text = ".localization-config" }
)

$wsDe.Range("A1").Hyperlinks.Delete()
foreach ($link in $deLinks) {
    $rng = $wsDe.Range($link.Cell)
    $wsDe.Hyperlinks.Add($rng, $link.Url, "", "", $link.Text)
    Set-LinkCellFormat $rng
}

Write-Host "Handoff report updated for b.md (zh-cn, de-de, Overview)"
